$valuesA = @(
[double]"-0.072771377589418762",
[double]"-0.061512152827962296",
[double]"-0.010498112877442622",
[double]"-0.0021322984136951106",
[double]"0.0011969251168091333",
[double]"-0.024162326551223856",
[double]"-0.014021216381922219",
[double]"-0.0040010000736860896",
[double]"-0.0019914745034683712",
[double]"1.1193389410024679e-05",
[double]"0.0030106310938844416",
[double]"0.0065103741060350551",
[double]"0.010019667537648047",
[double]"0.01803427616276565",
[double]"0.0038099774886068971",
[double]"0.0058272138631276782",
[double]"0.0078600448648273868",
[double]"-0.0048272611372155438",
[double]"-0.00075486532773583193",
[double]"-0.0080159702290014678",
[double]"-0.0040054454785973803",
[double]"-0.045714536408379303",
[double]"-0.040500718174895312",
[double]"-0.020099235629838397",
[double]"0.028475536339295004",
[double]"0.03103972233538066",
[double]"-0.041293722163110669",
[double]"-0.038754455639229413",
[double]"-0.031403101788436416",
[double]"0.028687123547392002",
[double]"0.035872145507815034",
[double]"-0.0040007695430652745"
)

$valuesB = @(
[double]"0.072339413916708395",
[double]"0.060212649248629901",
[double]"0.010132298344290191",
[double]"0.0018030748534130225",
[double]"-0.0023132098093139319",
[double]"0.024021216296825187",
[double]"0.014000999987194618",
[double]"0.0039914744800699786",
[double]"0.0019888065871285221",
[double]"-1.0631125272553277e-05",
[double]"-0.0030103741413132212",
[double]"-0.0065196675711969299",
[double]"-0.010034276231046135",
[double]"-0.018070145959115713",
[double]"-0.0038272138816939361",
[double]"-0.0058600448821453099",
[double]"-0.0078654463508174288",
[double]"0.0047548652949984671",
[double]"0.00022059351906555946",
[double]"0.0080054454450912971",
[double]"0.0039999999660853547",
[double]"0.045500718127042816",
[double]"0.040099235457343951",
[double]"0.019999999824698911",
[double]"-0.028539722359845143",
[double]"-0.031122452048844451",
[double]"0.040754455612963092",
[double]"0.038403101717379151",
[double]"0.031312875958612096",
[double]"-0.028872145576439134",
[double]"-0.035923654782545",
[double]"0.0039999999568038902"
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt 32; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $valuesA[$i]
    $ws.Cells.Item($row, 2).Value = $valuesB[$i]
}

# Column widths (stored OOXML "width" attribute goes through the host's
# character-width quantization, grid = 1/6; closest achievable snap to the
# target stored widths 16.42578125 / 16.28515625 are 16.5 / 16.333333...)
$ws.Columns.Item(1).ColumnWidth = 15.666666666666666
$ws.Columns.Item(2).ColumnWidth = 15.5
